$d = $word.ActiveDocument

function Set-ParagraphXml($paragraph, $innerBodyXml) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
          '<pkg:xmlData>' +
            '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:body>' + $innerBodyXml + '</w:body>' +
            '</w:document>' +
          '</pkg:xmlData>' +
        '</pkg:part>' +
      '</pkg:package>'
    $paragraph.Range.InsertXML($xml) | Out-Null
}

# 1. Rewrite the first body paragraph (intro / "I am excited to apply...") in place,
#    preserving the exact xml:space="preserve" run markup used throughout the document.
$intro = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "I am excited to apply*") {
        $intro = $p
        break
    }
}
$introXml = '<w:p><w:r><w:t xml:space="preserve">I am excited to apply for the Junior Software Engineer position at Compass. The role aligns perfectly with my skills and aspirations, especially in transforming education through innovative software solutions, a field that strongly interests me. Compass''s focus on collaborating on software development and deployment resonates with my passion - having developed a full-stack food ordering platform that increased international customer engagement by 10%, I understand the importance of effective communication in enhancing user satisfaction and platform efficiency. I am eager to contribute while growing with your team.</w:t></w:r></w:p>'
Set-ParagraphXml $intro $introXml

# 2. Locate the second body paragraph ("I am a Full stack Engineer...") and replace it
#    with its updated wording plus three new bulleted achievement items (separated by a
#    blank paragraph), reusing the document's existing bullet numbering definition
#    (numId 1 / ListParagraph style).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "I am a Full stack Engineer*") {
        $target = $p
        break
    }
}

$bodyXml =
    '<w:p><w:r><w:t xml:space="preserve">I am a Full stack Engineer who recently developed a full-stack food ordering platform for a caf' + [char]0xE9 + '. This experience strengthened my experience in React, Node.js, and MySQL and deepened my passion for solving practical challenges. A specific achievement from my previous experience that I believe can add value to the Junior Software Engineer position at Compass includes:</w:t></w:r></w:p>' +
    '<w:p/>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Increased international customer engagement by 10%.</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Proficient in secure user authentication and database design.</w:t></w:r></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Effective communication enhances user satisfaction and platform efficiency.</w:t></w:r></w:p>'

Set-ParagraphXml $target $bodyXml
